$wb = $excel.ActiveWorkbook

# ---- Home win ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'LEAGUE TWO'
$ws.Cells.Item(2,4).Value = 'Cheltenham - Salford City'
$ws.Cells.Item(2,5).Value = 73.3
$ws.Cells.Item(2,6).Value = 3
$ws.Cells.Item(3,1).Value = '25-01-2025 14:30'
$ws.Cells.Item(3,2).Value = 'GERMANY'
$ws.Cells.Item(3,3).Value = 'BUNDESLIGA'
$ws.Cells.Item(3,4).Value = 'FC Augsburg - 1. FC Heidenheim'
$ws.Cells.Item(3,5).Value = 73.3
$ws.Cells.Item(3,6).Value = 1.9
$ws.Cells.Item(4,1).Value = '25-01-2025 17:30'
$ws.Cells.Item(4,2).Value = 'GERMANY'
$ws.Cells.Item(4,3).Value = 'BUNDESLIGA'
$ws.Cells.Item(4,4).Value = 'Borussia Mönchengladbach - VfL Bochum'
$ws.Cells.Item(4,5).Value = 73.3
$ws.Cells.Item(4,6).Value = 1.7
$ws.Cells.Item(5,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(5,2).Value = 'GREECE'
$ws.Cells.Item(5,3).Value = 'SUPER LEAGUE 1'
$ws.Cells.Item(5,4).Value = 'OFI - Panserraikos'
$ws.Cells.Item(5,5).Value = 73.3
$ws.Cells.Item(5,6).Value = 1.91
$ws.Cells.Item(6,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(6,2).Value = 'SCOTLAND'
$ws.Cells.Item(6,3).Value = 'LEAGUE TWO'
$ws.Cells.Item(6,4).Value = 'Clyde - Elgin City'
$ws.Cells.Item(6,5).Value = 71.7
$ws.Cells.Item(6,6).Value = 2.37
$ws.Cells.Item(7,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(7,2).Value = 'ENGLAND'
$ws.Cells.Item(7,3).Value = 'NON LEAGUE PREMIER - NORTHERN'
$ws.Cells.Item(7,4).Value = 'Matlock Town - United Of Manchester'
$ws.Cells.Item(7,5).Value = 73.3
$ws.Cells.Item(7,6).Value = 2.2
$ws.Cells.Item(8,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(8,2).Value = 'ENGLAND'
$ws.Cells.Item(8,3).Value = 'NON LEAGUE PREMIER - SOUTHERN CENTRAL'
$ws.Cells.Item(8,4).Value = 'Harborough Town - Stourbridge'
$ws.Cells.Item(8,5).Value = 80
$ws.Cells.Item(8,6).Value = 2.2
$ws.Cells.Item(9,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(9,2).Value = 'ENGLAND'
$ws.Cells.Item(9,3).Value = 'NON LEAGUE PREMIER - SOUTHERN CENTRAL'
$ws.Cells.Item(9,4).Value = 'Stratford Town - Halesowen Town'
$ws.Cells.Item(9,5).Value = 73.3
$ws.Cells.Item(9,6).Value = 2.3
$ws.Cells.Item(10,1).Value = '25-01-2025 12:00'
$ws.Cells.Item(10,2).Value = 'GERMANY'
$ws.Cells.Item(10,3).Value = '2. BUNDESLIGA'
$ws.Cells.Item(10,4).Value = 'FC Schalke 04 - 1. FC Nürnberg'
$ws.Cells.Item(10,5).Value = 73.3
$ws.Cells.Item(10,6).Value = 2.1
$ws.Cells.Item(11,1).Value = '25-01-2025 00:00'
$ws.Cells.Item(11,2).Value = 'MEXICO'
$ws.Cells.Item(11,3).Value = 'LIGA PREMIER SERIE A'
$ws.Cells.Item(11,4).Value = 'Tecos - Real Apodaca'
$ws.Cells.Item(11,5).Value = 73.3
$ws.Cells.Item(11,6).Value = 2.1
$ws.Cells.Item(12,1).Value = '26-01-2025 17:00'
$ws.Cells.Item(12,2).Value = 'GUATEMALA'
$ws.Cells.Item(12,3).Value = 'LIGA NACIONAL'
$ws.Cells.Item(12,4).Value = 'Achuapa - Malacateco'
$ws.Cells.Item(12,5).Value = 73.3
$ws.Cells.Item(12,6).Value = 2.05
$ws.Cells.Item(13,1).Value = '26-01-2025 15:00'
$ws.Cells.Item(13,2).Value = 'PORTUGAL'
$ws.Cells.Item(13,3).Value = 'LIGA 3'
$ws.Cells.Item(13,4).Value = 'São João Ver - Varzim'
$ws.Cells.Item(13,5).Value = 73.3
$ws.Cells.Item(13,6).Value = 3
$ws.Cells.Item(14,1).Value = '26-01-2025 17:00'
$ws.Cells.Item(14,2).Value = 'SPAIN'
$ws.Cells.Item(14,3).Value = 'PRIMERA DIVISIÓN RFEF - GROUP 1'
$ws.Cells.Item(14,4).Value = 'Celta De Vigo II - Lugo'
$ws.Cells.Item(14,5).Value = 80
$ws.Cells.Item(14,6).Value = 1.8
$ws.Cells.Item(15,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(15,2).Value = 'SPAIN'
$ws.Cells.Item(15,3).Value = 'SEGUNDA DIVISIÓN RFEF - GROUP 1'
$ws.Cells.Item(15,4).Value = 'Compostela - Deportivo La Coruña II'
$ws.Cells.Item(15,5).Value = 73.3
$ws.Cells.Item(15,6).Value = 2.2
$ws.Cells.Item(16,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(16,2).Value = 'SPAIN'
$ws.Cells.Item(16,3).Value = 'SEGUNDA DIVISIÓN RFEF - GROUP 5'
$ws.Cells.Item(16,4).Value = 'Guadalajara - Tenerife II'
$ws.Cells.Item(16,5).Value = 80
$ws.Cells.Item(16,6).Value = 1.7
$ws.Cells.Item(17,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(17,2).Value = 'TURKEY'
$ws.Cells.Item(17,3).Value = '2. LIG'
$ws.Cells.Item(17,4).Value = 'Batman Petrolspor - Kastamonuspor 1966'
$ws.Cells.Item(17,5).Value = 70
$ws.Cells.Item(17,6).Value = 2

# ---- Draw ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = '25-01-2025 12:30'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'CHAMPIONSHIP'
$ws.Cells.Item(2,4).Value = 'Luton - Millwall'
$ws.Cells.Item(2,5).Value = 66.7
$ws.Cells.Item(2,6).Value = 3.5
$ws.Cells.Item(3,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(3,2).Value = 'ENGLAND'
$ws.Cells.Item(3,3).Value = 'LEAGUE TWO'
$ws.Cells.Item(3,4).Value = 'Fleetwood Town - Carlisle'
$ws.Cells.Item(3,5).Value = 66.7
$ws.Cells.Item(3,6).Value = 3.6
$ws.Cells.Item(4,1).Value = '25-01-2025 13:00'
$ws.Cells.Item(4,2).Value = 'FRANCE'
$ws.Cells.Item(4,3).Value = 'LIGUE 2'
$ws.Cells.Item(4,4).Value = 'Paris FC - RED Star FC 93'
$ws.Cells.Item(4,5).Value = 66.7
$ws.Cells.Item(4,6).Value = 4
$ws.Cells.Item(5,1).Value = '25-01-2025 14:00'
$ws.Cells.Item(5,2).Value = 'ITALY'
$ws.Cells.Item(5,3).Value = 'SERIE B'
$ws.Cells.Item(5,4).Value = 'Cesena - Bari'
$ws.Cells.Item(5,5).Value = 60
$ws.Cells.Item(5,6).Value = 3.25
$ws.Cells.Item(6,1).Value = '25-01-2025 18:30'
$ws.Cells.Item(6,2).Value = 'BRAZIL'
$ws.Cells.Item(6,3).Value = 'BRASILIENSE'
$ws.Cells.Item(6,4).Value = 'Ceilandense - Ceilândia'
$ws.Cells.Item(6,5).Value = 60
$ws.Cells.Item(6,6).Value = 7
$ws.Cells.Item(7,1).Value = '25-01-2025 18:00'
$ws.Cells.Item(7,2).Value = 'BRAZIL'
$ws.Cells.Item(7,3).Value = 'PAULISTA - A2'
$ws.Cells.Item(7,4).Value = 'São José EC - Primavera SP'
$ws.Cells.Item(7,5).Value = 60
$ws.Cells.Item(7,6).Value = 3.4
$ws.Cells.Item(8,1).Value = '25-01-2025 15:30'
$ws.Cells.Item(8,2).Value = 'GERMANY'
$ws.Cells.Item(8,3).Value = '3. LIGA'
$ws.Cells.Item(8,4).Value = 'SV Sandhausen - FC Saarbrücken'
$ws.Cells.Item(8,5).Value = 60
$ws.Cells.Item(8,6).Value = 3.35
$ws.Cells.Item(9,1).Value = '25-01-2025 13:00'
$ws.Cells.Item(9,2).Value = 'GERMANY'
$ws.Cells.Item(9,3).Value = '3. LIGA'
$ws.Cells.Item(9,4).Value = 'Verl - Alemannia Aachen'
$ws.Cells.Item(9,5).Value = 63.3
$ws.Cells.Item(9,6).Value = 3.5
$ws.Cells.Item(10,1).Value = '25-01-2025 12:30'
$ws.Cells.Item(10,2).Value = 'IRAN'
$ws.Cells.Item(10,3).Value = 'PERSIAN GULF PRO LEAGUE'
$ws.Cells.Item(10,4).Value = 'Esteghlal Khuzestan - Gol Gohar'
$ws.Cells.Item(10,5).Value = 66.7
$ws.Cells.Item(10,6).Value = 2.8
$ws.Cells.Item(11,1).Value = '25-01-2025 14:00'
$ws.Cells.Item(11,2).Value = 'ITALY'
$ws.Cells.Item(11,3).Value = 'SERIE C - GIRONE A'
$ws.Cells.Item(11,4).Value = 'Clodiense - Arzignano Valchiampo'
$ws.Cells.Item(11,5).Value = 73.3
$ws.Cells.Item(11,6).Value = 3.1
$ws.Cells.Item(12,1).Value = '25-01-2025 22:30'
$ws.Cells.Item(12,2).Value = 'PARAGUAY'
$ws.Cells.Item(12,3).Value = 'DIVISION PROFESIONAL - APERTURA'
$ws.Cells.Item(12,4).Value = 'General Caballero - Deportivo Recoleta'
$ws.Cells.Item(12,5).Value = 80
$ws.Cells.Item(12,6).Value = 3.3
$ws.Cells.Item(13,1).Value = '25-01-2025 13:10'
$ws.Cells.Item(13,2).Value = 'UNITED-ARAB-EMIRATES'
$ws.Cells.Item(13,3).Value = 'DIVISION 1'
$ws.Cells.Item(13,4).Value = 'Masfut - Al-Dhafra'
$ws.Cells.Item(13,5).Value = 60
$ws.Cells.Item(13,6).Value = 3.25
$ws.Cells.Item(14,1).Value = '26-01-2025 14:00'
$ws.Cells.Item(14,2).Value = 'ITALY'
$ws.Cells.Item(14,3).Value = 'SERIE B'
$ws.Cells.Item(14,4).Value = 'Brescia - Catanzaro'
$ws.Cells.Item(14,5).Value = 86.7
$ws.Cells.Item(14,6).Value = 3.1
$ws.Cells.Item(15,1).Value = '26-01-2025 12:30'
$ws.Cells.Item(15,2).Value = 'EGYPT'
$ws.Cells.Item(15,3).Value = 'SECOND LEAGUE'
$ws.Cells.Item(15,4).Value = 'Aswan Sc - Abu Qair Semad'
$ws.Cells.Item(15,5).Value = 60
$ws.Cells.Item(15,6).Value = 2.8
$ws.Cells.Item(16,1).Value = '26-01-2025 01:15'
$ws.Cells.Item(16,2).Value = 'EL-SALVADOR'
$ws.Cells.Item(16,3).Value = 'PRIMERA DIVISION'
$ws.Cells.Item(16,4).Value = 'Isidro Metapán - Once Municipal'
$ws.Cells.Item(16,5).Value = 60
$ws.Cells.Item(16,6).Value = 3.25
$ws.Cells.Item(17,1).Value = '26-01-2025 01:00'
$ws.Cells.Item(17,2).Value = 'GUATEMALA'
$ws.Cells.Item(17,3).Value = 'LIGA NACIONAL'
$ws.Cells.Item(17,4).Value = 'Marquense - Guastatoya'
$ws.Cells.Item(17,5).Value = 66.7
$ws.Cells.Item(17,6).Value = 3.5
$ws.Cells.Item(18,1).Value = '26-01-2025 15:30'
$ws.Cells.Item(18,2).Value = 'PORTUGAL'
$ws.Cells.Item(18,3).Value = 'SEGUNDA LIGA'
$ws.Cells.Item(18,4).Value = 'FC Porto B - Tondela'
$ws.Cells.Item(18,5).Value = 73.3
$ws.Cells.Item(18,6).Value = 3.1
$ws.Cells.Item(19,1).Value = '26-01-2025 17:30'
$ws.Cells.Item(19,2).Value = 'SPAIN'
$ws.Cells.Item(19,3).Value = 'SEGUNDA DIVISIÓN'
$ws.Cells.Item(19,4).Value = 'Racing Ferrol - Burgos'
$ws.Cells.Item(19,5).Value = 60
$ws.Cells.Item(19,6).Value = 2.8

# ---- Btts ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'CHAMPIONSHIP'
$ws.Cells.Item(2,4).Value = 'Preston - Middlesbrough'
$ws.Cells.Item(2,5).Value = 76.7
$ws.Cells.Item(2,6).Value = 1.73
$ws.Cells.Item(3,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(3,2).Value = 'ENGLAND'
$ws.Cells.Item(3,3).Value = 'LEAGUE ONE'
$ws.Cells.Item(3,4).Value = 'Exeter City - Blackpool'
$ws.Cells.Item(3,5).Value = 78.3
$ws.Cells.Item(3,6).Value = 1.75
$ws.Cells.Item(4,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(4,2).Value = 'ENGLAND'
$ws.Cells.Item(4,3).Value = 'LEAGUE ONE'
$ws.Cells.Item(4,4).Value = 'Stockport County - Crawley Town'
$ws.Cells.Item(4,5).Value = 83.3
$ws.Cells.Item(4,6).Value = 1.95
$ws.Cells.Item(5,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(5,2).Value = 'SCOTLAND'
$ws.Cells.Item(5,3).Value = 'PREMIERSHIP'
$ws.Cells.Item(5,4).Value = 'Ross County - Hibernian'
$ws.Cells.Item(5,5).Value = 76.7
$ws.Cells.Item(5,6).Value = 1.75
$ws.Cells.Item(6,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(6,2).Value = 'ENGLAND'
$ws.Cells.Item(6,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(6,4).Value = 'Ebbsfleet United - Braintree'
$ws.Cells.Item(6,5).Value = 82.2
$ws.Cells.Item(6,6).Value = 1.7
$ws.Cells.Item(7,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(7,2).Value = 'SCOTLAND'
$ws.Cells.Item(7,3).Value = 'LEAGUE TWO'
$ws.Cells.Item(7,4).Value = 'Forfar Athletic - Stranraer'
$ws.Cells.Item(7,5).Value = 86.7
$ws.Cells.Item(7,6).Value = 1.91
$ws.Cells.Item(8,1).Value = '25-01-2025 16:00'
$ws.Cells.Item(8,2).Value = 'BAHRAIN'
$ws.Cells.Item(8,3).Value = 'KING''S CUP'
$ws.Cells.Item(8,4).Value = 'Al-Hidd - Malkiya'
$ws.Cells.Item(8,5).Value = 76.7
$ws.Cells.Item(8,6).Value = 1.8
$ws.Cells.Item(9,1).Value = '25-01-2025 18:00'
$ws.Cells.Item(9,2).Value = 'BRAZIL'
$ws.Cells.Item(9,3).Value = 'POTIGUAR'
$ws.Cells.Item(9,4).Value = 'Santa Cruz RN - Força E Luz'
$ws.Cells.Item(9,5).Value = 90
$ws.Cells.Item(9,6).Value = 1.83
$ws.Cells.Item(10,1).Value = '25-01-2025 13:00'
$ws.Cells.Item(10,2).Value = 'GERMANY'
$ws.Cells.Item(10,3).Value = '3. LIGA'
$ws.Cells.Item(10,4).Value = 'Verl - Alemannia Aachen'
$ws.Cells.Item(10,5).Value = 75
$ws.Cells.Item(10,6).Value = 1.83
$ws.Cells.Item(11,1).Value = '25-01-2025 21:00'
$ws.Cells.Item(11,2).Value = 'HONDURAS'
$ws.Cells.Item(11,3).Value = 'LIGA NACIONAL'
$ws.Cells.Item(11,4).Value = 'CD Marathon - Victoria'
$ws.Cells.Item(11,5).Value = 80
$ws.Cells.Item(11,6).Value = 1.7
$ws.Cells.Item(12,1).Value = '25-01-2025 14:00'
$ws.Cells.Item(12,2).Value = 'ITALY'
$ws.Cells.Item(12,3).Value = 'SERIE C - GIRONE A'
$ws.Cells.Item(12,4).Value = 'Clodiense - Arzignano Valchiampo'
$ws.Cells.Item(12,5).Value = 86.7
$ws.Cells.Item(12,6).Value = 1.73
$ws.Cells.Item(13,1).Value = '25-01-2025 16:30'
$ws.Cells.Item(13,2).Value = 'ITALY'
$ws.Cells.Item(13,3).Value = 'SERIE C - GIRONE B'
$ws.Cells.Item(13,4).Value = 'Città Di Campobasso - Gubbio'
$ws.Cells.Item(13,5).Value = 76
$ws.Cells.Item(13,6).Value = 2.38
$ws.Cells.Item(14,1).Value = '25-01-2025 09:30'
$ws.Cells.Item(14,2).Value = 'MYANMAR'
$ws.Cells.Item(14,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(14,4).Value = 'Rakhine United - Yangon United'
$ws.Cells.Item(14,5).Value = 80.8
$ws.Cells.Item(14,6).Value = 1.85
$ws.Cells.Item(15,1).Value = '25-01-2025 11:00'
$ws.Cells.Item(15,2).Value = 'PORTUGAL'
$ws.Cells.Item(15,3).Value = 'SEGUNDA LIGA'
$ws.Cells.Item(15,4).Value = 'Felgueiras 1932 - Alverca'
$ws.Cells.Item(15,5).Value = 83.3
$ws.Cells.Item(15,6).Value = 1.8
$ws.Cells.Item(16,1).Value = '25-01-2025 14:00'
$ws.Cells.Item(16,2).Value = 'PORTUGAL'
$ws.Cells.Item(16,3).Value = 'SEGUNDA LIGA'
$ws.Cells.Item(16,4).Value = 'Leixoes - Chaves'
$ws.Cells.Item(16,5).Value = 75
$ws.Cells.Item(16,6).Value = 1.95
$ws.Cells.Item(17,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(17,2).Value = 'SPAIN'
$ws.Cells.Item(17,3).Value = 'PRIMERA DIVISIÓN RFEF - GROUP 2'
$ws.Cells.Item(17,4).Value = 'Recreativo Huelva - AD Ceuta FC'
$ws.Cells.Item(17,5).Value = 76.7
$ws.Cells.Item(17,6).Value = 1.95
$ws.Cells.Item(18,1).Value = '26-01-2025 16:30'
$ws.Cells.Item(18,2).Value = 'ENGLAND'
$ws.Cells.Item(18,3).Value = 'PREMIER LEAGUE'
$ws.Cells.Item(18,4).Value = 'Aston Villa - West Ham'
$ws.Cells.Item(18,5).Value = 80
$ws.Cells.Item(18,6).Value = 1.75
$ws.Cells.Item(19,1).Value = '26-01-2025 11:30'
$ws.Cells.Item(19,2).Value = 'ITALY'
$ws.Cells.Item(19,3).Value = 'SERIE A'
$ws.Cells.Item(19,4).Value = 'AC Milan - Parma'
$ws.Cells.Item(19,5).Value = 76.7
$ws.Cells.Item(19,6).Value = 1.75
$ws.Cells.Item(20,1).Value = '26-01-2025 13:30'
$ws.Cells.Item(20,2).Value = 'NETHERLANDS'
$ws.Cells.Item(20,3).Value = 'EREDIVISIE'
$ws.Cells.Item(20,4).Value = 'Waalwijk - Willem II'
$ws.Cells.Item(20,5).Value = 77.8
$ws.Cells.Item(20,6).Value = 1.75
$ws.Cells.Item(21,1).Value = '26-01-2025 14:00'
$ws.Cells.Item(21,2).Value = 'ITALY'
$ws.Cells.Item(21,3).Value = 'SERIE B'
$ws.Cells.Item(21,4).Value = 'Brescia - Catanzaro'
$ws.Cells.Item(21,5).Value = 90
$ws.Cells.Item(21,6).Value = 1.8
$ws.Cells.Item(22,1).Value = '26-01-2025 17:30'
$ws.Cells.Item(22,2).Value = 'BELGIUM'
$ws.Cells.Item(22,3).Value = 'JUPILER PRO LEAGUE'
$ws.Cells.Item(22,4).Value = 'Anderlecht - KV Mechelen'
$ws.Cells.Item(22,5).Value = 76.7
$ws.Cells.Item(22,6).Value = 1.7
$ws.Cells.Item(23,1).Value = '26-01-2025 13:00'
$ws.Cells.Item(23,2).Value = 'BRAZIL'
$ws.Cells.Item(23,3).Value = 'MINEIRO - 1'
$ws.Cells.Item(23,4).Value = 'Uberlandia - Itabirito'
$ws.Cells.Item(23,5).Value = 86.7
$ws.Cells.Item(23,6).Value = 1.85
$ws.Cells.Item(24,1).Value = '26-01-2025 21:30'
$ws.Cells.Item(24,2).Value = 'BRAZIL'
$ws.Cells.Item(24,3).Value = 'PARANAENSE - 1'
$ws.Cells.Item(24,4).Value = 'Londrina - Parana'
$ws.Cells.Item(24,5).Value = 75.8
$ws.Cells.Item(24,6).Value = 1.91
$ws.Cells.Item(25,1).Value = '26-01-2025 21:00'
$ws.Cells.Item(25,2).Value = 'CHILE'
$ws.Cells.Item(25,3).Value = 'COPA CHILE'
$ws.Cells.Item(25,4).Value = 'U. Catolica - Everton De Vina'
$ws.Cells.Item(25,5).Value = 80
$ws.Cells.Item(25,6).Value = 1.83
$ws.Cells.Item(26,1).Value = '26-01-2025 13:30'
$ws.Cells.Item(26,2).Value = 'ITALY'
$ws.Cells.Item(26,3).Value = 'SERIE D - GIRONE C'
$ws.Cells.Item(26,4).Value = 'Mestre - Caravaggio'
$ws.Cells.Item(26,5).Value = 80
$ws.Cells.Item(26,6).Value = 2.4
$ws.Cells.Item(27,1).Value = '26-01-2025 01:05'
$ws.Cells.Item(27,2).Value = 'MEXICO'
$ws.Cells.Item(27,3).Value = 'LIGA DE EXPANSIÓN MX'
$ws.Cells.Item(27,4).Value = 'Cancún - Alebrijes De Oaxaca'
$ws.Cells.Item(27,5).Value = 80.8
$ws.Cells.Item(27,6).Value = 1.8
$ws.Cells.Item(28,1).Value = '26-01-2025 23:45'
$ws.Cells.Item(28,2).Value = 'PARAGUAY'
$ws.Cells.Item(28,3).Value = 'DIVISION PROFESIONAL - APERTURA'
$ws.Cells.Item(28,4).Value = 'Cerro Porteno - Libertad Asuncion'
$ws.Cells.Item(28,5).Value = 83.3
$ws.Cells.Item(28,6).Value = 1.93
$ws.Cells.Item(29,1).Value = '26-01-2025 15:30'
$ws.Cells.Item(29,2).Value = 'PORTUGAL'
$ws.Cells.Item(29,3).Value = 'SEGUNDA LIGA'
$ws.Cells.Item(29,4).Value = 'FC Porto B - Tondela'
$ws.Cells.Item(29,5).Value = 90
$ws.Cells.Item(29,6).Value = 1.8
$ws.Cells.Item(30,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(30,2).Value = 'PORTUGAL'
$ws.Cells.Item(30,3).Value = 'SEGUNDA LIGA'
$ws.Cells.Item(30,4).Value = 'Penafiel - Vizela'
$ws.Cells.Item(30,5).Value = 81.1
$ws.Cells.Item(30,6).Value = 1.8
$ws.Cells.Item(31,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(31,2).Value = 'TURKEY'
$ws.Cells.Item(31,3).Value = '2. LIG'
$ws.Cells.Item(31,4).Value = 'Batman Petrolspor - Kastamonuspor 1966'
$ws.Cells.Item(31,5).Value = 76
$ws.Cells.Item(31,6).Value = 1.85

# ---- Over_Under ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'CHAMPIONSHIP'
$ws.Cells.Item(2,4).Value = 'QPR - Sheffield Wednesday'
$ws.Cells.Item(2,5).Value = 80
$ws.Cells.Item(2,6).Value = 1.8
$ws.Cells.Item(2,7).Value = 30
$ws.Cells.Item(2,8).Value = 3.2
$ws.Cells.Item(3,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(3,2).Value = 'ENGLAND'
$ws.Cells.Item(3,3).Value = 'LEAGUE ONE'
$ws.Cells.Item(3,4).Value = 'Stockport County - Crawley Town'
$ws.Cells.Item(3,5).Value = 80
$ws.Cells.Item(3,6).Value = 1.8
$ws.Cells.Item(3,7).Value = 32.5
$ws.Cells.Item(3,8).Value = 3
$ws.Cells.Item(4,1).Value = '25-01-2025 14:30'
$ws.Cells.Item(4,2).Value = 'GERMANY'
$ws.Cells.Item(4,3).Value = 'BUNDESLIGA'
$ws.Cells.Item(4,4).Value = 'Borussia Dortmund - Werder Bremen'
$ws.Cells.Item(4,5).Value = 60
$ws.Cells.Item(4,6).Value = 1.5
$ws.Cells.Item(4,7).Value = 60
$ws.Cells.Item(4,8).Value = 2.2
$ws.Cells.Item(5,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(5,2).Value = 'SCOTLAND'
$ws.Cells.Item(5,3).Value = 'PREMIERSHIP'
$ws.Cells.Item(5,4).Value = 'Ross County - Hibernian'
$ws.Cells.Item(5,5).Value = 80
$ws.Cells.Item(5,6).Value = 1.91
$ws.Cells.Item(5,7).Value = 40
$ws.Cells.Item(5,8).Value = 3.2
$ws.Cells.Item(6,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(6,2).Value = 'ENGLAND'
$ws.Cells.Item(6,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(6,4).Value = 'AFC Fylde - Wealdstone'
$ws.Cells.Item(6,5).Value = 85
$ws.Cells.Item(6,6).Value = 1.7
$ws.Cells.Item(6,7).Value = 50
$ws.Cells.Item(6,8).Value = 2.62
$ws.Cells.Item(7,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(7,2).Value = 'ENGLAND'
$ws.Cells.Item(7,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(7,4).Value = 'Aldershot Town - Tamworth'
$ws.Cells.Item(7,5).Value = 80
$ws.Cells.Item(7,6).Value = 1.7
$ws.Cells.Item(7,7).Value = 53.3
$ws.Cells.Item(7,8).Value = 2.62
$ws.Cells.Item(8,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(8,2).Value = 'ENGLAND'
$ws.Cells.Item(8,3).Value = 'NATIONAL LEAGUE'
$ws.Cells.Item(8,4).Value = 'Gateshead - Oldham'
$ws.Cells.Item(8,5).Value = 90
$ws.Cells.Item(8,6).Value = 1.73
$ws.Cells.Item(8,7).Value = 67.5
$ws.Cells.Item(8,8).Value = 2.75
$ws.Cells.Item(9,1).Value = '25-01-2025 16:00'
$ws.Cells.Item(9,2).Value = 'BAHRAIN'
$ws.Cells.Item(9,3).Value = 'KING''S CUP'
$ws.Cells.Item(9,4).Value = 'Al Ittifaq Maqaba - Isa Town'
$ws.Cells.Item(9,5).Value = 77.5
$ws.Cells.Item(9,6).Value = 1.85
$ws.Cells.Item(9,7).Value = 67.5
$ws.Cells.Item(9,8).Value = 3
$ws.Cells.Item(10,1).Value = '25-01-2025 19:00'
$ws.Cells.Item(10,2).Value = 'BELGIUM'
$ws.Cells.Item(10,3).Value = 'CHALLENGER PRO LEAGUE'
$ws.Cells.Item(10,4).Value = 'RSC Anderlecht II - Patro Eisden'
$ws.Cells.Item(10,5).Value = 70
$ws.Cells.Item(10,6).Value = 1.67
$ws.Cells.Item(10,7).Value = 60
$ws.Cells.Item(10,8).Value = 2.7
$ws.Cells.Item(11,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(11,2).Value = 'BELGIUM'
$ws.Cells.Item(11,3).Value = 'CHALLENGER PRO LEAGUE'
$ws.Cells.Item(11,4).Value = 'Seraing United - Francs Borains'
$ws.Cells.Item(11,5).Value = 75
$ws.Cells.Item(11,6).Value = 1.73
$ws.Cells.Item(11,7).Value = 60
$ws.Cells.Item(11,8).Value = 2.7
$ws.Cells.Item(12,1).Value = '25-01-2025 19:00'
$ws.Cells.Item(12,2).Value = 'BELGIUM'
$ws.Cells.Item(12,3).Value = 'FIRST AMATEUR DIVISION'
$ws.Cells.Item(12,4).Value = 'Cappellen - Hoogstraten'
$ws.Cells.Item(12,5).Value = 80
$ws.Cells.Item(12,6).Value = 1.5
$ws.Cells.Item(12,7).Value = 61.8
$ws.Cells.Item(12,8).Value = 2.25
$ws.Cells.Item(13,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(13,2).Value = 'CYPRUS'
$ws.Cells.Item(13,3).Value = '1. DIVISION'
$ws.Cells.Item(13,4).Value = 'Omonia 29is Maiou - Apollon Limassol'
$ws.Cells.Item(13,5).Value = 85
$ws.Cells.Item(13,6).Value = 1.91
$ws.Cells.Item(13,7).Value = 45
$ws.Cells.Item(13,8).Value = 3.5
$ws.Cells.Item(14,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(14,2).Value = 'ENGLAND'
$ws.Cells.Item(14,3).Value = 'NATIONAL LEAGUE - NORTH'
$ws.Cells.Item(14,4).Value = 'Oxford City - Kidderminster Harriers'
$ws.Cells.Item(14,5).Value = 90
$ws.Cells.Item(14,6).Value = 2.1
$ws.Cells.Item(14,7).Value = 55
$ws.Cells.Item(14,8).Value = 3.8
$ws.Cells.Item(15,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(15,2).Value = 'ENGLAND'
$ws.Cells.Item(15,3).Value = 'NATIONAL LEAGUE - SOUTH'
$ws.Cells.Item(15,4).Value = 'Weston-super-Mare - Welling United'
$ws.Cells.Item(15,5).Value = 60
$ws.Cells.Item(15,6).Value = 1.85
$ws.Cells.Item(15,7).Value = 60
$ws.Cells.Item(15,8).Value = 3.3
$ws.Cells.Item(16,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(16,2).Value = 'ENGLAND'
$ws.Cells.Item(16,3).Value = 'NON LEAGUE PREMIER - SOUTHERN CENTRAL'
$ws.Cells.Item(16,4).Value = 'Bedford Town - Banbury United'
$ws.Cells.Item(16,5).Value = 80
$ws.Cells.Item(16,6).Value = 1.85
$ws.Cells.Item(16,7).Value = 60
$ws.Cells.Item(16,8).Value = 3.3
$ws.Cells.Item(17,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(17,2).Value = 'ENGLAND'
$ws.Cells.Item(17,3).Value = 'NON LEAGUE PREMIER - SOUTHERN SOUTH'
$ws.Cells.Item(17,4).Value = 'Dorchester Town - Sholing'
$ws.Cells.Item(17,5).Value = 80
$ws.Cells.Item(17,6).Value = 1.75
$ws.Cells.Item(17,7).Value = 40
$ws.Cells.Item(17,8).Value = 2.88
$ws.Cells.Item(18,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(18,2).Value = 'ENGLAND'
$ws.Cells.Item(18,3).Value = 'NON LEAGUE PREMIER - SOUTHERN SOUTH'
$ws.Cells.Item(18,4).Value = 'Frome Town - Gloucester City'
$ws.Cells.Item(18,5).Value = 73.3
$ws.Cells.Item(18,6).Value = 1.6
$ws.Cells.Item(18,7).Value = 60
$ws.Cells.Item(18,8).Value = 2.55
$ws.Cells.Item(19,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(19,2).Value = 'ENGLAND'
$ws.Cells.Item(19,3).Value = 'NON LEAGUE PREMIER - SOUTHERN SOUTH'
$ws.Cells.Item(19,4).Value = 'Gosport Borough - Walton & Hersham'
$ws.Cells.Item(19,5).Value = 70
$ws.Cells.Item(19,6).Value = 1.5
$ws.Cells.Item(19,7).Value = 70
$ws.Cells.Item(19,8).Value = 2.25
$ws.Cells.Item(20,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(20,2).Value = 'ENGLAND'
$ws.Cells.Item(20,3).Value = 'NON LEAGUE PREMIER - SOUTHERN SOUTH'
$ws.Cells.Item(20,4).Value = 'Winchester City - Swindon Supermarine'
$ws.Cells.Item(20,5).Value = 90
$ws.Cells.Item(20,6).Value = 1.45
$ws.Cells.Item(20,7).Value = 80
$ws.Cells.Item(20,8).Value = 2.2
$ws.Cells.Item(21,1).Value = '25-01-2025 15:30'
$ws.Cells.Item(21,2).Value = 'GERMANY'
$ws.Cells.Item(21,3).Value = '3. LIGA'
$ws.Cells.Item(21,4).Value = 'SV Sandhausen - FC Saarbrücken'
$ws.Cells.Item(21,5).Value = 70
$ws.Cells.Item(21,6).Value = 1.8
$ws.Cells.Item(21,7).Value = 60
$ws.Cells.Item(21,8).Value = 3
$ws.Cells.Item(22,1).Value = '25-01-2025 15:30'
$ws.Cells.Item(22,2).Value = 'MALTA'
$ws.Cells.Item(22,3).Value = 'CHALLENGE LEAGUE'
$ws.Cells.Item(22,4).Value = 'Mgarr United - Lija Athletic'
$ws.Cells.Item(22,5).Value = 80
$ws.Cells.Item(22,6).Value = 1.61
$ws.Cells.Item(22,7).Value = 60
$ws.Cells.Item(22,8).Value = 2.5
$ws.Cells.Item(23,1).Value = '25-01-2025 18:00'
$ws.Cells.Item(23,2).Value = 'MALTA'
$ws.Cells.Item(23,3).Value = 'PREMIER LEAGUE'
$ws.Cells.Item(23,4).Value = 'Zabbar St. Patrick - Balzan FC'
$ws.Cells.Item(23,5).Value = 80
$ws.Cells.Item(23,6).Value = 2
$ws.Cells.Item(23,7).Value = 46.7
$ws.Cells.Item(23,8).Value = 3.4
$ws.Cells.Item(24,1).Value = '25-01-2025 00:00'
$ws.Cells.Item(24,2).Value = 'MEXICO'
$ws.Cells.Item(24,3).Value = 'LIGA PREMIER SERIE A'
$ws.Cells.Item(24,4).Value = 'Tecos - Real Apodaca'
$ws.Cells.Item(24,5).Value = 80
$ws.Cells.Item(24,6).Value = 1.7
$ws.Cells.Item(24,7).Value = 35
$ws.Cells.Item(24,8).Value = 2.7
$ws.Cells.Item(25,1).Value = '25-01-2025 14:30'
$ws.Cells.Item(25,2).Value = 'NETHERLANDS'
$ws.Cells.Item(25,3).Value = 'TWEEDE DIVISIE'
$ws.Cells.Item(25,4).Value = 'Katwijk - Koninklijke HFC'
$ws.Cells.Item(25,5).Value = 65
$ws.Cells.Item(25,6).Value = 1.7
$ws.Cells.Item(25,7).Value = 60
$ws.Cells.Item(25,8).Value = 2.7
$ws.Cells.Item(26,1).Value = '25-01-2025 12:50'
$ws.Cells.Item(26,2).Value = 'SAUDI-ARABIA'
$ws.Cells.Item(26,3).Value = 'DIVISION 1'
$ws.Cells.Item(26,4).Value = 'Al Taee - Al Suqoor'
$ws.Cells.Item(26,5).Value = 80
$ws.Cells.Item(26,6).Value = 1.95
$ws.Cells.Item(26,7).Value = 26.7
$ws.Cells.Item(26,8).Value = 3.3
$ws.Cells.Item(27,1).Value = '25-01-2025 15:15'
$ws.Cells.Item(27,2).Value = 'SPAIN'
$ws.Cells.Item(27,3).Value = 'SEGUNDA DIVISIÓN'
$ws.Cells.Item(27,4).Value = 'Cordoba - Racing Santander'
$ws.Cells.Item(27,5).Value = 93.3
$ws.Cells.Item(27,6).Value = 1.77
$ws.Cells.Item(27,7).Value = 33.3
$ws.Cells.Item(27,8).Value = 2.88
$ws.Cells.Item(28,1).Value = '25-01-2025 19:30'
$ws.Cells.Item(28,2).Value = 'SWITZERLAND'
$ws.Cells.Item(28,3).Value = 'SUPER LEAGUE'
$ws.Cells.Item(28,4).Value = 'FC Winterthur - FC Lugano'
$ws.Cells.Item(28,5).Value = 88.8
$ws.Cells.Item(28,6).Value = 1.6
$ws.Cells.Item(28,7).Value = 68.8
$ws.Cells.Item(28,8).Value = 2.45
$ws.Cells.Item(29,1).Value = '25-01-2025 11:30'
$ws.Cells.Item(29,2).Value = 'THAILAND'
$ws.Cells.Item(29,3).Value = 'THAI LEAGUE 2'
$ws.Cells.Item(29,4).Value = 'Lampang FC - Chonburi FC'
$ws.Cells.Item(29,5).Value = 75
$ws.Cells.Item(29,6).Value = 1.75
$ws.Cells.Item(29,7).Value = 60
$ws.Cells.Item(29,8).Value = 2.9
$ws.Cells.Item(30,1).Value = '25-01-2025 12:00'
$ws.Cells.Item(30,2).Value = 'THAILAND'
$ws.Cells.Item(30,3).Value = 'THAI LEAGUE 2'
$ws.Cells.Item(30,4).Value = 'Suphanburi - Chanthaburi'
$ws.Cells.Item(30,5).Value = 85
$ws.Cells.Item(30,6).Value = 1.73
$ws.Cells.Item(30,7).Value = 35
$ws.Cells.Item(30,8).Value = 2.75
$ws.Cells.Item(31,1).Value = '25-01-2025 10:30'
$ws.Cells.Item(31,2).Value = 'TURKEY'
$ws.Cells.Item(31,3).Value = '1. LIG'
$ws.Cells.Item(31,4).Value = 'Boluspor - Genclerbirligi'
$ws.Cells.Item(31,5).Value = 73.8
$ws.Cells.Item(31,6).Value = 2.15
$ws.Cells.Item(31,7).Value = 61.3
$ws.Cells.Item(31,8).Value = 3.8
$ws.Cells.Item(32,1).Value = '25-01-2025 13:00'
$ws.Cells.Item(32,2).Value = 'TURKEY'
$ws.Cells.Item(32,3).Value = '1. LIG'
$ws.Cells.Item(32,4).Value = 'İstanbulspor - Amed'
$ws.Cells.Item(32,5).Value = 60
$ws.Cells.Item(32,6).Value = 2
$ws.Cells.Item(32,7).Value = 60
$ws.Cells.Item(32,8).Value = 3.4
$ws.Cells.Item(33,1).Value = '25-01-2025 23:30'
$ws.Cells.Item(33,2).Value = 'WORLD'
$ws.Cells.Item(33,3).Value = 'SUDAMERICANO U20'
$ws.Cells.Item(33,4).Value = 'Peru U20 - Venezuela U20'
$ws.Cells.Item(33,5).Value = 80
$ws.Cells.Item(33,6).Value = 2.05
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(33,8).Value = 3.65
$ws.Cells.Item(34,1).Value = '26-01-2025 19:45'
$ws.Cells.Item(34,2).Value = 'FRANCE'
$ws.Cells.Item(34,3).Value = 'LIGUE 1'
$ws.Cells.Item(34,4).Value = 'Nice - Marseille'
$ws.Cells.Item(34,5).Value = 85
$ws.Cells.Item(34,6).Value = 1.73
$ws.Cells.Item(34,7).Value = 45
$ws.Cells.Item(34,8).Value = 2.75
$ws.Cells.Item(35,1).Value = '26-01-2025 17:30'
$ws.Cells.Item(35,2).Value = 'GREECE'
$ws.Cells.Item(35,3).Value = 'SUPER LEAGUE 1'
$ws.Cells.Item(35,4).Value = 'PAOK - Levadiakos'
$ws.Cells.Item(35,5).Value = 80
$ws.Cells.Item(35,6).Value = 1.7
$ws.Cells.Item(35,7).Value = 48.8
$ws.Cells.Item(35,8).Value = 2.62
$ws.Cells.Item(36,1).Value = '26-01-2025 06:00'
$ws.Cells.Item(36,2).Value = 'AUSTRALIA'
$ws.Cells.Item(36,3).Value = 'A-LEAGUE'
$ws.Cells.Item(36,4).Value = 'Western Sydney Wanderers - Auckland'
$ws.Cells.Item(36,5).Value = 85
$ws.Cells.Item(36,6).Value = 1.73
$ws.Cells.Item(36,7).Value = 85
$ws.Cells.Item(36,8).Value = 2.75
$ws.Cells.Item(37,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(37,2).Value = 'SPAIN'
$ws.Cells.Item(37,3).Value = 'SEGUNDA DIVISIÓN RFEF - GROUP 2'
$ws.Cells.Item(37,4).Value = 'Real Zaragoza II - Izarra'
$ws.Cells.Item(37,5).Value = 90
$ws.Cells.Item(37,6).Value = 2.15
$ws.Cells.Item(37,7).Value = 75
$ws.Cells.Item(37,8).Value = 3.8
$ws.Cells.Item(38,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(38,2).Value = 'SPAIN'
$ws.Cells.Item(38,3).Value = 'SEGUNDA DIVISIÓN RFEF - GROUP 5'
$ws.Cells.Item(38,4).Value = 'Colonia Moscardó - SS Reyes'
$ws.Cells.Item(38,5).Value = 73.3
$ws.Cells.Item(38,6).Value = 2.25
$ws.Cells.Item(38,7).Value = 60
$ws.Cells.Item(38,8).Value = 4
$ws.Cells.Item(39,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(39,2).Value = 'TURKEY'
$ws.Cells.Item(39,3).Value = '2. LIG'
$ws.Cells.Item(39,4).Value = 'Batman Petrolspor - Kastamonuspor 1966'
$ws.Cells.Item(39,5).Value = 80
$ws.Cells.Item(39,6).Value = 2.15
$ws.Cells.Item(39,7).Value = 60
$ws.Cells.Item(39,8).Value = 3.8
$ws.Cells.Item(40,1).Value = '26-01-2025 11:00'
$ws.Cells.Item(40,2).Value = 'TURKEY'
$ws.Cells.Item(40,3).Value = '2. LIG'
$ws.Cells.Item(40,4).Value = 'Van BB - Nazilli Belediyespor'
$ws.Cells.Item(40,5).Value = 80
$ws.Cells.Item(40,6).Value = 1.73
$ws.Cells.Item(40,7).Value = 45
$ws.Cells.Item(40,8).Value = 2.88

# ---- Away Win ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2,1).Value = '25-01-2025 15:00'
$ws.Cells.Item(2,2).Value = 'ENGLAND'
$ws.Cells.Item(2,3).Value = 'NON LEAGUE PREMIER - ISTHMIAN'
$ws.Cells.Item(2,4).Value = 'Dulwich Hamlet - Chichester City'
$ws.Cells.Item(2,5).Value = 70
$ws.Cells.Item(2,6).Value = 2.45
$ws.Cells.Item(3,1).Value = '25-01-2025 13:00'
$ws.Cells.Item(3,2).Value = 'WORLD'
$ws.Cells.Item(3,3).Value = 'FRIENDLIES CLUBS'
$ws.Cells.Item(3,4).Value = 'Sogndal - Hønefoss'
$ws.Cells.Item(3,5).Value = 90
$ws.Cells.Item(3,6).Value = 7
$ws.Cells.Item(4,1).Value = '26-01-2025 21:30'
$ws.Cells.Item(4,2).Value = 'EL-SALVADOR'
$ws.Cells.Item(4,3).Value = 'PRIMERA DIVISION'
$ws.Cells.Item(4,4).Value = 'Fuerte San Francisco - Municipal Limeño'
$ws.Cells.Item(4,5).Value = 73.3
$ws.Cells.Item(4,6).Value = 2.6

Write-Host "done"
